# regenerate orders with updates distance/sizes
#
# The stimulus set was regenerated with new distance codes and a new "far"
# size code:
#   D64 -> D69
#   D80 -> D86
#   D51 -> D55
#   S30 -> S31   (S20 / S25 are unchanged)
#
# Every cell whose text contains one of the old tokens (Condition,
# Filename_Left, Filename_Right, Distance, Size columns) needs the token(s)
# swapped for the new one. Numbers / booleans / the Face / ConditionID
# columns are untouched since none of their values embed these tokens.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow  = $usedRange.Row
$firstCol  = $usedRange.Column
$lastRow   = $firstRow + $usedRange.Rows.Count - 1
$lastCol   = $firstCol + $usedRange.Columns.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $old = $cell.Value()
        if ($old -is [string]) {
            $new = $old.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
            if ($new -ne $old) {
                $cell.Value = $new
            }
        }
    }
}
